$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 14170.167
$ws.Range("J34").Value = 14610.2
$ws.Range("L34").Value = 14610.2
$ws.Range("N34").Value = -15016.2
$ws.Range("H36").Value = 14170.167
$ws.Range("J36").Value = 14610.2
$ws.Range("L36").Value = 14610.2
$ws.Range("N36").Value = -16040.2
$ws.Range("H54").Value = 14374.625
$ws.Range("I54").Value = 10000
$ws.Range("K54").Value = 10000
$ws.Range("M54").Value = -9514
$ws.Range("H63").Value = 30000.166
$ws.Range("J63").Value = 30000.166
$ws.Range("L63").Value = 30000.166
$ws.Range("N63").Value = -31248.166
$ws.Range("H66").Value = 30000.166
$ws.Range("J66").Value = 30000.166
$ws.Range("L66").Value = 90000.49800000001
$ws.Range("N66").Value = -96240.49800000001
$ws.Range("H100").Value = 2330
$ws.Range("I100").Value = 1602
$ws.Range("J100").Value = 3058
$ws.Range("K100").Value = 1602
$ws.Range("L100").Value = 3058
$ws.Range("M100").Value = -1061
$ws.Range("N100").Value = -4140
$ws.Range("H141").Value = 7245.6895
$ws.Range("I141").Value = 4005.4
$ws.Range("K141").Value = 12016.2
$ws.Range("M141").Value = -6836.200000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 158112.72
$ws.Range("I102").Value = 167790.33
$ws.Range("J102").Value = 100047
$ws.Range("K102").Value = 167790.33
$ws.Range("L102").Value = 100047
$ws.Range("M102").Value = -166168.33
$ws.Range("N102").Value = -103291

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 273
$ws.Range("I64").Value = 200
$ws.Range("K64").Value = 200
$ws.Range("M64").Value = 25
$ws.Range("H67").Value = 273
$ws.Range("I67").Value = 200
$ws.Range("K67").Value = 200
$ws.Range("M67").Value = 580
$ws.Range("H80").Value = 546.1818
$ws.Range("I80").Value = 522.55554
$ws.Range("J80").Value = 562.53845
$ws.Range("K80").Value = 522.55554
$ws.Range("L80").Value = 562.53845
$ws.Range("M80").Value = 475.44446
$ws.Range("N80").Value = -2558.53845
$ws.Range("H83").Value = 546.1818
$ws.Range("I83").Value = 522.55554
$ws.Range("J83").Value = 562.53845
$ws.Range("K83").Value = 2612.7777
$ws.Range("L83").Value = 2812.69225
$ws.Range("M83").Value = 2379.2223
$ws.Range("N83").Value = -12796.69225
$ws.Range("H134").Value = 1198.3208
$ws.Range("I134").Value = 903.6905
$ws.Range("J134").Value = 2323.2727
$ws.Range("K134").Value = 2711.0715
$ws.Range("L134").Value = 6969.8181
$ws.Range("M134").Value = -176.0715
$ws.Range("N134").Value = -12039.8181

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1712.15
$ws.Range("I31").Value = 1011.5536
$ws.Range("J31").Value = 2603.818
$ws.Range("K31").Value = 1011.5536
$ws.Range("L31").Value = 2603.818
$ws.Range("M31").Value = -716.5536
$ws.Range("N31").Value = -3193.818
$ws.Range("H34").Value = 1712.15
$ws.Range("I34").Value = 1011.5536
$ws.Range("J34").Value = 2603.818
$ws.Range("K34").Value = 1011.5536
$ws.Range("L34").Value = 2603.818
$ws.Range("M34").Value = -809.5536
$ws.Range("N34").Value = -3007.818
$ws.Range("H58").Value = 29412870
$ws.Range("I58").Value = 41667560
$ws.Range("J58").Value = 1615.8
$ws.Range("K58").Value = 41667560
$ws.Range("L58").Value = 1615.8
$ws.Range("M58").Value = -41667357
$ws.Range("N58").Value = -2021.8
$ws.Range("H136").Value = 29412870
$ws.Range("I136").Value = 41667560
$ws.Range("J136").Value = 1615.8
$ws.Range("K136").Value = 125002680
$ws.Range("L136").Value = 4847.4
$ws.Range("M136").Value = -125000130
$ws.Range("N136").Value = -9947.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 857.2593000000001
$ws.Range("I5").Value = 589.86664
$ws.Range("J5").Value = 1191.5
$ws.Range("K5").Value = 1769.59992
$ws.Range("L5").Value = 3574.5
$ws.Range("M5").Value = -1657.59992
$ws.Range("N5").Value = -3798.5
$ws.Range("H13").Value = 87.75
$ws.Range("I13").Value = 83.666664
$ws.Range("J13").Value = 100
$ws.Range("K13").Value = 250.999992
$ws.Range("L13").Value = 300
$ws.Range("M13").Value = -82.99999199999999
$ws.Range("N13").Value = -636
$ws.Range("H22").Value = 200503.9
$ws.Range("I22").Value = 666767
$ws.Range("J22").Value = 676.8570999999999
$ws.Range("K22").Value = 2000301
$ws.Range("L22").Value = 2030.5713
$ws.Range("M22").Value = -2000132
$ws.Range("N22").Value = -2368.5713
$ws.Range("H27").Value = 200503.9
$ws.Range("I27").Value = 666767
$ws.Range("J27").Value = 676.8570999999999
$ws.Range("K27").Value = 2000301
$ws.Range("L27").Value = 2030.5713
$ws.Range("M27").Value = -2000199
$ws.Range("N27").Value = -2234.5713
$ws.Range("H60").Value = 500
$ws.Range("I60").Value = 250
$ws.Range("J60").Value = 1000
$ws.Range("K60").Value = 750
$ws.Range("L60").Value = 3000
$ws.Range("M60").Value = -499
$ws.Range("N60").Value = -3502
$ws.Range("H122").Value = 705.8929000000001
$ws.Range("I122").Value = 304.6842
$ws.Range("J122").Value = 1552.8889
$ws.Range("K122").Value = 2742.1578
$ws.Range("L122").Value = 13976.0001
$ws.Range("M122").Value = -292.1578
$ws.Range("N122").Value = -18876.0001
$ws.Range("H135").Value = 857.2593000000001
$ws.Range("I135").Value = 589.86664
$ws.Range("J135").Value = 1191.5
$ws.Range("K135").Value = 5308.79976
$ws.Range("L135").Value = 10723.5
$ws.Range("M135").Value = -2773.79976
$ws.Range("N135").Value = -15793.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H63").Value = 10000
$ws.Range("J63").Value = 10000
$ws.Range("L63").Value = 10000
$ws.Range("N63").Value = -11372
$ws.Range("H66").Value = 10000
$ws.Range("J66").Value = 10000
$ws.Range("L66").Value = 30000
$ws.Range("N66").Value = -36864
$ws.Range("H92").Value = 5237
$ws.Range("J92").Value = 6750.5
$ws.Range("L92").Value = 6750.5
$ws.Range("N92").Value = -10494.5
$ws.Range("H134").Value = 24000
$ws.Range("J134").Value = 24000
$ws.Range("L134").Value = 72000
$ws.Range("N134").Value = -77070

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H105").Value = 12880
$ws.Range("J105").Value = 12880
$ws.Range("L105").Value = 12880
$ws.Range("N105").Value = -19868
$ws.Range("H139").Value = 55000
$ws.Range("J139").Value = 55000
$ws.Range("L139").Value = 55000
$ws.Range("N139").Value = -65280
$ws.Range("H141").Value = 160000
$ws.Range("J141").Value = 160000
$ws.Range("L141").Value = 160000
$ws.Range("N141").Value = -170360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 11882.125
$ws.Range("J64").Value = 11882.125
$ws.Range("L64").Value = 11882.125
$ws.Range("N64").Value = -12378.125
$ws.Range("H67").Value = 11882.125
$ws.Range("J67").Value = 11882.125
$ws.Range("L67").Value = 11882.125
$ws.Range("N67").Value = -13598.125
$ws.Range("H132").Value = 6013.567
$ws.Range("I132").Value = 8412.111000000001
$ws.Range("J132").Value = 2415.75
$ws.Range("K132").Value = 25236.333
$ws.Range("L132").Value = 7247.25
$ws.Range("M132").Value = -22706.333
$ws.Range("N132").Value = -12307.25
